$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.868.74"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "1.667.08"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.530"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.11%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.253"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0896"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.03%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "1.677.47"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "26.876.44"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -1.98%  "
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").Value = "1.463.72"
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.571"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0168"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.977"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.43%  "
$ws.Range("D45").Value = "1.812.60"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("E51").Value = "  +1.15%  "
